# Generate Report for Handback
# Refresh the handback-status report's "latest generate / handoff / handback"
# timestamps for the f79a32b9-... file across the Overview, zh-cn and de-de
# sheets, as produced by a fresh run of the handback report generator.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-17 12:45:39"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-17 12:45:34"
$zhcn.Range("K3").Value = "2016-08-17 12:45:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-17 12:45:39"
$dede.Range("K3").Value = "2016-08-17 12:46:06"
